$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.035.18"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.876.36"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.37"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -3.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2927"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06599"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("D10").Value = "1.882.69"
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.62"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07176"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6673"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.22"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.923"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "29.970.07"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007792"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9987"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("D20").Value = "2.118.89"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9991"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.772"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.868"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.098"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.79"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.67"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.93"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.894"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.382"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.188"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08748"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.975"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05014"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7174"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.111"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.661"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01816"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.685"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.156"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9299"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.761"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9981"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4217"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.372"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.24%  "
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.75"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3763"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.247"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.339"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.16%  "
